$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 202.14815
$ws.Range("I9").Value = 206.07692
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 206.07692
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -37.07692
$ws.Range("N9").Value = -438

# Row 21
$ws.Range("H21").Value = 36618
$ws.Range("J21").Value = 29937.8
$ws.Range("L21").Value = 29937.8
$ws.Range("N21").Value = -30873.8

# Row 23
$ws.Range("H23").Value = 36618
$ws.Range("J23").Value = 29937.8
$ws.Range("L23").Value = 29937.8
$ws.Range("N23").Value = -30405.8

# Row 137
$ws.Range("H137").Value = 2625.9119
$ws.Range("I137").Value = 1473
$ws.Range("K137").Value = 4419
$ws.Range("M137").Value = -1869

# Row 138
$ws.Range("H138").Value = 2783.78
$ws.Range("I138").Value = 630.7273
$ws.Range("J138").Value = 3844.2388
$ws.Range("K138").Value = 1892.1819
$ws.Range("L138").Value = 11532.7164
$ws.Range("M138").Value = 3247.8181
$ws.Range("N138").Value = -21812.7164

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1262.5
$ws.Range("I2").Value = 1305.6666
$ws.Range("J2").Value = 1133
$ws.Range("K2").Value = 1305.6666
$ws.Range("L2").Value = 1133
$ws.Range("M2").Value = -1192.6666
$ws.Range("N2").Value = -1359

# Row 32
$ws.Range("H32").Value = 5416.0967
$ws.Range("I32").Value = 4266.857
$ws.Range("K32").Value = 4266.857
$ws.Range("M32").Value = -3979.857

# Row 88
$ws.Range("H88").Value = 7411151.5
$ws.Range("I88").Value = 9527481
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 9527481
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -9527075
$ws.Range("N88").Value = -4812

# Row 91
$ws.Range("H91").Value = 7411151.5
$ws.Range("I91").Value = 9527481
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 9527481
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -9526077
$ws.Range("N91").Value = -6808

# Row 110
$ws.Range("H110").Value = 1106.1052
$ws.Range("I110").Value = 1140.1666
$ws.Range("J110").Value = 493
$ws.Range("K110").Value = 1140.1666
$ws.Range("L110").Value = 493
$ws.Range("M110").Value = 904.8334
$ws.Range("N110").Value = -4583

# Row 116
$ws.Range("H116").Value = 1262.5
$ws.Range("I116").Value = 1305.6666
$ws.Range("J116").Value = 1133
$ws.Range("K116").Value = 1305.6666
$ws.Range("L116").Value = 1133
$ws.Range("M116").Value = 988.3334
$ws.Range("N116").Value = -5721

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1262.5
$ws.Range("I3").Value = 1305.6666
$ws.Range("J3").Value = 1133
$ws.Range("K3").Value = 1305.6666
$ws.Range("L3").Value = 1133
$ws.Range("M3").Value = -1191.6666
$ws.Range("N3").Value = -1361

# Row 86
$ws.Range("H86").Value = 2054.3103
$ws.Range("I86").Value = 1843.0526
$ws.Range("J86").Value = 2455.7
$ws.Range("K86").Value = 1843.0526
$ws.Range("L86").Value = 2455.7
$ws.Range("M86").Value = -720.0526
$ws.Range("N86").Value = -4701.7

# Row 89
$ws.Range("H89").Value = 2054.3103
$ws.Range("I89").Value = 1843.0526
$ws.Range("J89").Value = 2455.7
$ws.Range("K89").Value = 9215.262999999999
$ws.Range("L89").Value = 12278.5
$ws.Range("M89").Value = -3599.262999999999
$ws.Range("N89").Value = -23510.5

# Row 95
$ws.Range("H95").Value = 33500
$ws.Range("J95").Value = 33500
$ws.Range("L95").Value = 33500
$ws.Range("N95").Value = -38992

# Row 103
$ws.Range("H103").Value = 24513.305
$ws.Range("J103").Value = 24513.305
$ws.Range("L103").Value = 24513.305
$ws.Range("N103").Value = -26857.305

# Row 134
$ws.Range("H134").Value = 1719
$ws.Range("I134").Value = 1273.5106
$ws.Range("J134").Value = 6953.5
$ws.Range("K134").Value = 3820.5318
$ws.Range("L134").Value = 20860.5
$ws.Range("M134").Value = -1285.5318
$ws.Range("N134").Value = -25930.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 6537177.5
$ws.Range("I16").Value = 10102047
$ws.Range("J16").Value = 1583
$ws.Range("K16").Value = 10102047
$ws.Range("L16").Value = 1583
$ws.Range("M16").Value = -10101760
$ws.Range("N16").Value = -2157

# Row 31
$ws.Range("H31").Value = 11630183
$ws.Range("I31").Value = 1139.6296
$ws.Range("K31").Value = 1139.6296
$ws.Range("M31").Value = -844.6296

# Row 34
$ws.Range("H34").Value = 11630183
$ws.Range("I34").Value = 1139.6296
$ws.Range("K34").Value = 1139.6296
$ws.Range("M34").Value = -937.6296

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 105
$ws.Range("H105").Value = 2452.8
$ws.Range("I105").Value = 2075.4285
$ws.Range("K105").Value = 2075.4285
$ws.Range("M105").Value = -328.4285

# Row 113
$ws.Range("H113").Value = 6537177.5
$ws.Range("I113").Value = 10102047
$ws.Range("J113").Value = 1583
$ws.Range("K113").Value = 10102047
$ws.Range("L113").Value = 1583
$ws.Range("M113").Value = -10099877
$ws.Range("N113").Value = -5923

# Row 132
$ws.Range("H132").Value = 3486.625
$ws.Range("I132").Value = 2377.4
$ws.Range("K132").Value = 7132.200000000001
$ws.Range("M132").Value = -4602.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 75
$ws.Range("H75").Value = 7010.75
$ws.Range("I75").Value = 1013
$ws.Range("J75").Value = 9010
$ws.Range("K75").Value = 3039
$ws.Range("L75").Value = 27030
$ws.Range("M75").Value = -2041
$ws.Range("N75").Value = -29026

# Row 78
$ws.Range("H78").Value = 7010.75
$ws.Range("I78").Value = 1013
$ws.Range("J78").Value = 9010
$ws.Range("K78").Value = 9117
$ws.Range("L78").Value = 81090
$ws.Range("M78").Value = -4125
$ws.Range("N78").Value = -91074

# Row 131
$ws.Range("H131").Value = 8197537
$ws.Range("I131").Value = 55555800
$ws.Range("J131").Value = 914.5769
$ws.Range("K131").Value = 166667400
$ws.Range("L131").Value = 2743.7307
$ws.Range("M131").Value = -166662360
$ws.Range("N131").Value = -12823.7307

# Row 133
$ws.Range("H133").Value = 3093
$ws.Range("J133").Value = 4166.6665
$ws.Range("L133").Value = 12499.9995
$ws.Range("N133").Value = -22619.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 6192.7617
$ws.Range("I70").Value = 5823.906
$ws.Range("J70").Value = 8147.7
$ws.Range("K70").Value = 5823.906
$ws.Range("L70").Value = 8147.7
$ws.Range("M70").Value = -5553.906
$ws.Range("N70").Value = -8687.700000000001

# Row 73
$ws.Range("H73").Value = 6192.7617
$ws.Range("I73").Value = 5823.906
$ws.Range("J73").Value = 8147.7
$ws.Range("K73").Value = 5823.906
$ws.Range("L73").Value = 8147.7
$ws.Range("M73").Value = -4887.906
$ws.Range("N73").Value = -10019.7

# Row 113
$ws.Range("H113").Value = 1656.0555
$ws.Range("I113").Value = 1789.0769
$ws.Range("K113").Value = 1789.0769
$ws.Range("M113").Value = 380.9231

# Row 126
$ws.Range("H126").Value = 3289.3
$ws.Range("I126").Value = 2950.375
$ws.Range("J126").Value = 4645
$ws.Range("K126").Value = 8851.125
$ws.Range("L126").Value = 13935
$ws.Range("M126").Value = -6381.125
$ws.Range("N126").Value = -18875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 26
$ws.Range("H26").Value = 29987
$ws.Range("J26").Value = 29987
$ws.Range("L26").Value = 29987
$ws.Range("N26").Value = -30577

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Row 42
$ws.Range("H42").Value = 34797.5
$ws.Range("J42").Value = 34797.5
$ws.Range("L42").Value = 34797.5
$ws.Range("N42").Value = -35923.5

# Row 49
$ws.Range("H49").Value = 34797.5
$ws.Range("J49").Value = 34797.5
$ws.Range("L49").Value = 34797.5
$ws.Range("N49").Value = -35091.5

# Row 132
$ws.Range("H132").Value = 5131.913
$ws.Range("I132").Value = 1496
$ws.Range("J132").Value = 13442.571
$ws.Range("K132").Value = 4488
$ws.Range("L132").Value = 40327.713
$ws.Range("M132").Value = -1958
$ws.Range("N132").Value = -45387.713

# Row 136
$ws.Range("H136").Value = 3407.0715
$ws.Range("I136").Value = 1406.7222
$ws.Range("J136").Value = 7007.7
$ws.Range("K136").Value = 4220.1666
$ws.Range("L136").Value = 21023.1
$ws.Range("M136").Value = -1670.1666
$ws.Range("N136").Value = -26123.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 54
$ws.Range("H54").Value = 17894.572
$ws.Range("J54").Value = 17894.572
$ws.Range("L54").Value = 17894.572
$ws.Range("N54").Value = -18934.572

# Row 62
$ws.Range("H62").Value = 27500
$ws.Range("I62").Value = 3333.3333
$ws.Range("K62").Value = 3333.3333
$ws.Range("M62").Value = -2709.3333

# Row 65
$ws.Range("H65").Value = 27500
$ws.Range("I65").Value = 3333.3333
$ws.Range("K65").Value = 16666.6665
$ws.Range("M65").Value = -13546.6665

# Row 108
$ws.Range("H108").Value = 39750
$ws.Range("J108").Value = 39750
$ws.Range("L108").Value = 39750
$ws.Range("N108").Value = -47430
